$wb = $excel.ActiveWorkbook

# --- Features sheet: insert two new columns (M, N) before the current
#     last column ("remarks"), add headers + comments, and size the
#     new columns. The existing M column ("remarks") shifts to O.
$ws = $wb.Worksheets.Item("Features")
$ws.Activate()

$ws.Columns("M:N").Insert()

$ws.Range("M1").Value = "Curve_Fit_Model"
$ws.Range("N1").Value = "Curve_Fit_Weighting"

$ws.Columns("M:N").ColumnWidth = 19.83

$ws.Range("M1").AddComment("Either 'linear' or 'quadratic'") | Out-Null
$ws.Range("N1").AddComment("Either '1/x', '1/x2' or 'none'") | Out-Null

$ws.Range("O11").Select() | Out-Null

# --- Info sheet becomes the active/selected tab
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Activate()

Write-Host "done"
